$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New MG Firing - Advancing Fire event goes in right after "e054a" (row 80)
# and before "e060" (current row 81), so insert a fresh row at 81 and push
# everything from the old row 81 down by one.
$ws.Rows(81).Insert()

# New row gets its own height (ht="180") and the e054b content.
$ws.Rows(81).RowHeight = 180

$ws.Cells.Item(81, 1).Value2 = "e054b"

$body = "<Bold>e054b MG Firing - Advancing Fire</Bold> " + "`r`n" + `
"<InlineUIContainer><Button Content='r4.74.3' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   " + "`r`n" + `
"<InlineUIContainer><Button Content='r22.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> " + "`r`n" + `
"<LineBreak/><LineBreak/>" + "`r`n" + `
"Advancing fire attacks infantry targets that move into the zone. It also help protect against Panzerfaust attacks. " + "`r`n" + `
"<LineBreak/><LineBreak/>" + "`r`n" + `
"Roll 2D for ammo expected:<LineBreak/>" + "`r`n" + `
"01 - 30 = 1 box expended<LineBreak/>" + "`r`n" + `
"31 - 97 = no effect<LineBreak/>" + "`r`n" + `
"98 - 100 = MG malfunctions" + "`r`n" + `
"<LineBreak/><LineBreak/>" + "`r`n" + `
"Die Roll = <InlineUIContainer><Image Name='DieRollBlue' Height='21' Width='21' > </Image></InlineUIContainer>"

$ws.Cells.Item(81, 2).Value2 = $body

# Match the author's final scroll position / selection.
$excel.ActiveWindow.ScrollRow = 79
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B81").Select()
